$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 950
$ws.Range("I6").Value = 950
$ws.Range("K6").Value = 2850
$ws.Range("M6").Value = -2738

$ws.Range("H8").Value = 9186.625
$ws.Range("I8").Value = 582.1667
$ws.Range("J8").Value = 35000
$ws.Range("K8").Value = 1746.5001
$ws.Range("L8").Value = 105000
$ws.Range("M8").Value = -1607.5001
$ws.Range("N8").Value = -105278

$ws.Range("H31").Value = 2528824.5
$ws.Range("I31").Value = 2528824.5
$ws.Range("K31").Value = 7586473.5
$ws.Range("M31").Value = -7586243.5

$ws.Range("H40").Value = 1569.125
$ws.Range("I40").Value = 1361.3572
$ws.Range("J40").Value = 1860
$ws.Range("K40").Value = 1361.3572
$ws.Range("L40").Value = 1860
$ws.Range("M40").Value = -1186.3572
$ws.Range("N40").Value = -2210

$ws.Range("H58").Value = 1366.0667
$ws.Range("I58").Value = 192.33333
$ws.Range("J58").Value = 2539.8
$ws.Range("K58").Value = 576.99999
$ws.Range("L58").Value = 7619.400000000001
$ws.Range("M58").Value = -426.99999
$ws.Range("N58").Value = -7919.400000000001

$ws.Range("H97").Value = 1356
$ws.Range("J97").Value = 1410.7693
$ws.Range("L97").Value = 4232.3079
$ws.Range("N97").Value = -5224.3079

$ws.Range("H133").Value = 61068
$ws.Range("J133").Value = 61068
$ws.Range("L133").Value = 61068
$ws.Range("N133").Value = -71188

$ws.Range("H139").Value = 53398
$ws.Range("J139").Value = 53398
$ws.Range("L139").Value = 53398
$ws.Range("N139").Value = -63678

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2525.25
$ws.Range("I122").Value = 2105.5
$ws.Range("J122").Value = 2945
$ws.Range("K122").Value = 6316.5
$ws.Range("L122").Value = 8835
$ws.Range("M122").Value = -3866.5
$ws.Range("N122").Value = -13735

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 39975
$ws.Range("J35").Value = 39975
$ws.Range("L35").Value = 39975
$ws.Range("N35").Value = -40595

$ws.Range("H107").Value = 2381.923
$ws.Range("I107").Value = 2524.182
$ws.Range("J107").Value = 1599.5
$ws.Range("K107").Value = 2524.182
$ws.Range("L107").Value = 1599.5
$ws.Range("M107").Value = -604.1819999999998
$ws.Range("N107").Value = -5439.5

$ws.Range("H132").Value = 62679.824
$ws.Range("J132").Value = 62679.824
$ws.Range("L132").Value = 62679.824
$ws.Range("N132").Value = -72799.82399999999

$ws.Range("H138").Value = 57913.75
$ws.Range("J138").Value = 57913.75
$ws.Range("L138").Value = 57913.75
$ws.Range("N138").Value = -68193.75

$ws.Range("H140").Value = 48394.375
$ws.Range("J140").Value = 48394.375
$ws.Range("L140").Value = 48394.375
$ws.Range("N140").Value = -58754.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 27875
$ws.Range("I17").Value = 5750
$ws.Range("J17").Value = 50000
$ws.Range("K17").Value = 5750
$ws.Range("L17").Value = 50000
$ws.Range("M17").Value = -5576
$ws.Range("N17").Value = -50348

$ws.Range("H97").Value = 34992.8
$ws.Range("J97").Value = 34992.8
$ws.Range("L97").Value = 34992.8
$ws.Range("N97").Value = -36974.8

$ws.Range("H135").Value = 39475.715
$ws.Range("J135").Value = 39475.715
$ws.Range("L135").Value = 39475.715
$ws.Range("N135").Value = -49615.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1609.7084
$ws.Range("I5").Value = 1918.6
$ws.Range("J5").Value = 1094.8889
$ws.Range("K5").Value = 5755.799999999999
$ws.Range("L5").Value = 3284.6667
$ws.Range("M5").Value = -5643.799999999999
$ws.Range("N5").Value = -3508.6667

$ws.Range("H64").Value = 5195.375
$ws.Range("I64").Value = 3603.75
$ws.Range("K64").Value = 10811.25
$ws.Range("M64").Value = -10541.25

$ws.Range("H67").Value = 5195.375
$ws.Range("I67").Value = 3603.75
$ws.Range("K67").Value = 10811.25
$ws.Range("M67").Value = -9875.25

$ws.Range("H122").Value = 652.9268
$ws.Range("J122").Value = 623.2121
$ws.Range("L122").Value = 5608.908899999999
$ws.Range("N122").Value = -10508.9089

$ws.Range("H131").Value = 16950658
$ws.Range("J131").Value = 18520116
$ws.Range("L131").Value = 55560348
$ws.Range("N131").Value = -55570428

$ws.Range("H135").Value = 1609.7084
$ws.Range("I135").Value = 1918.6
$ws.Range("J135").Value = 1094.8889
$ws.Range("K135").Value = 17267.4
$ws.Range("L135").Value = 9854.000099999999
$ws.Range("M135").Value = -14732.4
$ws.Range("N135").Value = -14924.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = ""

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = ""

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""

$ws.Range("H126").Value = 2284.25
$ws.Range("I126").Value = 2401.8333
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 7205.499899999999
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -4735.499899999999
$ws.Range("N126").Value = -11440.0001

$ws.Range("H132").Value = 3750
$ws.Range("I132").Value = 3283.6667
$ws.Range("K132").Value = 9851.000100000001
$ws.Range("M132").Value = -7321.000100000001

$ws.Range("H140").Value = 85881.664
$ws.Range("J140").Value = 85881.664
$ws.Range("L140").Value = 85881.664
$ws.Range("N140").Value = -96241.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 717.2222
$ws.Range("I22").Value = 563.75
$ws.Range("J22").Value = 840
$ws.Range("K22").Value = 563.75
$ws.Range("L22").Value = 840
$ws.Range("M22").Value = -268.75
$ws.Range("N22").Value = -1430

$ws.Range("H27").Value = 717.2222
$ws.Range("I27").Value = 563.75
$ws.Range("J27").Value = 840
$ws.Range("K27").Value = 563.75
$ws.Range("L27").Value = 840
$ws.Range("M27").Value = -456.75
$ws.Range("N27").Value = -1054

$ws.Range("H46").Value = 1480
$ws.Range("I46").Value = 1350
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1350
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1162
$ws.Range("N46").Value = -2376

$ws.Range("H136").Value = 3308.8572
$ws.Range("I136").Value = 2697.524
$ws.Range("K136").Value = 8092.572
$ws.Range("M136").Value = -5542.572

$ws.Range("H137").Value = 69765.8
$ws.Range("I137").Value = 49000
$ws.Range("J137").Value = 83609.664
$ws.Range("K137").Value = 49000
$ws.Range("L137").Value = 83609.664
$ws.Range("M137").Value = -43900
$ws.Range("N137").Value = -93809.664

$ws.Range("H140").Value = 39792.5
$ws.Range("J140").Value = 39792.5
$ws.Range("L140").Value = 39792.5
$ws.Range("N140").Value = -50152.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -39920

$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws.Range("H141").Value = 70271.125
$ws.Range("J141").Value = 70271.125
$ws.Range("L141").Value = 70271.125
$ws.Range("N141").Value = -80631.125
